$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.689.30'
$ws.Range("E2").Value = '  -7.03%  '

$ws.Range("D3").Value = '2.544.77'
$ws.Range("E3").Value = '  -4.45%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '299.75'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.97%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '92.40'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.17%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.575'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -4.08%  '

$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.546'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -6.06%  '

$ws.Range("E10").Value = '  -6.51%  '

$ws.Range("E11").Value = '  -5.33%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.66'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -5.54%  '

$ws.Range("E13").Value = '  +5.30%  '

$ws.Range("D14").Value = '2.934.52'
$ws.Range("E14").Value = '  -4.04%  '

$ws.Range("D15").Value = '2.531.97'
$ws.Range("E15").Value = '  -4.56%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.875'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -5.95%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.21'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -5.80%  '

$ws.Range("D18").Value = '42.672.14'
$ws.Range("E18").Value = '  -7.05%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.87'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.37%  '

$ws.Range("D20").Value = '0.0₃0983'
$ws.Range("E20").Value = '  -3.68%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.58'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.50%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '71.62'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.54%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '255.46'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -9.46%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.91'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.86%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.12'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -5.95%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '29.08'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -6.17%  '

$ws.Range("E27").Value = '  -0.02%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.07'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.56%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '37.13'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.02%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.11'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.66%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.99'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.04%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '152.28'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.50%  '

$ws.Range("E33").Value = '  -2.24%  '

$ws.Range("E34").Value = '  -8.73%  '

$ws.Range("E35").Value = '  -9.48%  '

$ws.Range("E36").Value = '  -5.84%  '

$ws.Range("E37").Value = '  -5.65%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '24.12'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -6.03%  '

$ws.Range("E39").Value = '  -4.61%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '17.05'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.93%  '

$ws.Range("B41").Value = 'NEARProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.41'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.65%  '

$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0311'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.33%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.86'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.86%  '

$ws.Range("D44").Value = '2.083.12'
$ws.Range("E44").Value = '  -3.59%  '

$ws.Range("E45").Value = '  -0.03%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.63'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.71%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '84.53'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -10.12%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.01'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.33%  '

$ws.Range("D49").Value = '2.789.68'
$ws.Range("E49").Value = '  -4.11%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '104.47'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -6.81%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.67'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.94%  '
